$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.853.80"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.342.33"
$ws.Range("E3").Value = "  +6.33%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'258.30"
$ws.Range("E5").Value = "  +6.83%  "
$ws.Range("D6").Value = "'625.14"
$ws.Range("E6").Value = "  +2.61%  "
$ws.Range("D7").Value = "'1.43"
$ws.Range("E7").Value = "  +28.14%  "
$ws.Range("D8").Value = "'0.391"
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'0.877"
$ws.Range("E10").Value = "  +11.42%  "
$ws.Range("D11").Value = "3.341.52"
$ws.Range("E11").Value = "  +6.36%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "'37.38"
$ws.Range("E13").Value = "  +10.23%  "
$ws.Range("D14").Value = "98.416.26"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'0.0000249"
$ws.Range("E15").Value = "  +3.89%  "
$ws.Range("D16").Value = "3.951.02"
$ws.Range("E16").Value = "  +5.92%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "3.334.97"
$ws.Range("E18").Value = "  +6.26%  "
$ws.Range("D19").Value = "'3.55"
$ws.Range("E19").Value = "  +3.72%  "
$ws.Range("D20").Value = "'15.20"
$ws.Range("E20").Value = "  +4.70%  "
$ws.Range("D21").Value = "'491.25"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("E22").Value = "  +7.76%  "
$ws.Range("D23").Value = "'0.0000211"
$ws.Range("E23").Value = "  +10.27%  "
$ws.Range("D24").Value = "'9.36"
$ws.Range("E24").Value = "  +6.68%  "
$ws.Range("D25").Value = "'5.63"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").Value = "'88.90"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'11.85"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "3.499.06"
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D29").Value = "'0.293"
$ws.Range("E29").Value = "  +23.39%  "
$ws.Range("E31").Value = "  +11.82%  "
$ws.Range("D32").Value = "'0.139"
$ws.Range("E32").Value = "  +13.41%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.68"
$ws.Range("E33").Value = "  +8.61%  "
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'28.10"
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("D36").Value = "'0.150"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").Value = "'7.26"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("D39").Value = "'498.78"
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.460"
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'24.86"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("D43").Value = "'3.65"
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = "  +6.26%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'0.782"
$ws.Range("E45").Value = "  +12.57%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'159.86"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "'1.94"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").Value = "'0.849"
$ws.Range("E49").Value = "  +7.64%  "
$ws.Range("D50").Value = "'4.63"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").Value = "'45.69"
$ws.Range("E51").Value = "  +3.71%  "
